# Apply updated symbol list values (price/volume) to match the latest scrape.
# Cells are stored as text, so we force text entry with a leading apostrophe
# and then reset the style to Normal to avoid leaving a quote-prefix style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''297.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''-0.97%'
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = '''0.12%'
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''-0.97%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.07999'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''9.08%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''2.491'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''38.09%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''7.780'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''-0.06%'
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.9239'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''-0.05%'
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.1731'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''2.92%'
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.07364'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''3.44%'
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.08901'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''9.71%'
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.03033'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''0.42%'
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''0.86%'
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''0.001495'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''0.23%'
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.005951'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''-3.99%'
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''3.518'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''1.73%'
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''3.799'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''1.66%'
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''1.24%'
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.3254'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''0.87%'
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''1.52%'
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''4.288'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''-5.76%'
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''2.30%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.04596'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''-1.06%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.001242'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''2.30%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.004425'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''-6.68%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.0001200'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''-7.53%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.0003426'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''82.97%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = '''0.01771'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''3.35%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.04473'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''-0.51%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.006843'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''-3.39%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.1341'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''0.28%'
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.002209'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''-0.72%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.009800'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-6.13%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.00006561'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''5.15%'
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''-0.03%'
$ws.Range("E46").Style = "Normal"
$ws.Range("E48").Value = '''-57.27%'
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.00002100'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''-0.03%'
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.0002000'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''0.04%'
$ws.Range("E50").Style = "Normal"
